$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 145 - this shifts the existing rows
# 145..175 down to 146..176 (carrying their values/formatting along).
$ws.Rows.Item(145).Insert()

# Populate the newly inserted row 145 with the new weekly record.
$ws.Cells.Item(145, 1).Value = 11
$ws.Cells.Item(145, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(145, 3).Value = "Bíobío"
$ws.Cells.Item(145, 4).Value = 45127
$ws.Cells.Item(145, 5).Value = 8
$ws.Cells.Item(145, 6).Value = 100112001
$ws.Cells.Item(145, 7).Value = "Berenjena"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 160
$ws.Cells.Item(145, 11).Value = 8000
$ws.Cells.Item(145, 12).Value = 9000
$ws.Cells.Item(145, 13).Value = 8375
$ws.Cells.Item(145, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(145, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(145, 16).Value = 168
$ws.Cells.Item(145, 17).Value = 50
$ws.Cells.Item(145, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date number-format style
# used by the rest of the "Fecha" column (style index 2 in the
# original workbook == the YYYY-MM-DD HH:MM:SS date format).
$ws.Cells.Item(145, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
